$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the counts to reflect the two-variable (Fe-number, Frost) clustering
$ws.Range("B2").Value = 384
$ws.Range("B3").Value = 331

# Remove the now-obsolete third cluster row (area1 vs rest distinction only needs 2 rows)
$ws.Rows.Item(4).Delete()
